$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => updated Price (D) / Volume(1h) (E) values, per the commit diff.
# A leading apostrophe on D-column values forces Excel to keep them as
# text (matching the original inlineStr cell type) instead of silently
# parsing look-alike numbers (e.g. "246.48") into numeric cells.
$updates = @(
    @{ Row = 2; D = "'42.123.69"; E = "  -0.85%  " }
    @{ Row = 3; D = "'2.236.96"; E = "  -1.76%  " }
    @{ Row = 4; D = $null; E = "  +0.02%  " }
    @{ Row = 5; D = "'246.48"; E = "  -1.49%  " }
    @{ Row = 6; D = "'0.629"; E = "  -0.50%  " }
    @{ Row = 7; D = "'76.28"; E = "  +5.82%  " }
    @{ Row = 8; D = $null; E = "  +0.02%  " }
    @{ Row = 9; D = "'0.625"; E = "  -2.43%  " }
    @{ Row = 10; D = "'41.08"; E = "  +6.01%  " }
    @{ Row = 11; D = "'0.0945"; E = "  -1.28%  " }
    @{ Row = 12; D = "'7.12"; E = "  -2.23%  " }
    @{ Row = 13; D = $null; E = "  -1.12%  " }
    @{ Row = 14; D = "'2.573.35"; E = "  -1.70%  " }
    @{ Row = 15; D = "'14.78"; E = "  -1.76%  " }
    @{ Row = 16; D = "'0.857"; E = "  -2.13%  " }
    @{ Row = 17; D = "'2.235.70"; E = "  -2.20%  " }
    @{ Row = 18; D = "'41.991.31"; E = "  -1.06%  " }
    @{ Row = 19; D = "'0.0₃0977"; E = "  -1.63%  " }
    @{ Row = 20; D = "'6.12"; E = "  -2.46%  " }
    @{ Row = 21; D = "'71.35"; E = "  -1.00%  " }
    @{ Row = 22; D = "'230.24"; E = "  -0.87%  " }
    @{ Row = 23; D = "'2.19"; E = "  -2.28%  " }
    @{ Row = 24; D = $null; E = "  -0.02%  " }
    @{ Row = 25; D = "'3.70"; E = "  -5.67%  " }
    @{ Row = 26; D = "'11.14"; E = "  -3.04%  " }
    @{ Row = 27; D = "'2.32"; E = "  -4.23%  " }
    @{ Row = 28; D = "'7.27"; E = "  +13.82%  " }
    @{ Row = 29; D = $null; E = "  -1.69%  " }
    @{ Row = 30; D = "'169.15"; E = "  +1.32%  " }
    @{ Row = 31; D = "'20.47"; E = "  -2.60%  " }
    @{ Row = 32; D = "'0.0854"; E = "  +5.52%  " }
    @{ Row = 33; D = "'32.99"; E = "  +5.99%  " }
    @{ Row = 34; D = $null; E = "  -5.13%  " }
    @{ Row = 35; D = $null; E = "  +0.83%  " }
    @{ Row = 36; D = "'4.59"; E = "  -2.50%  " }
    @{ Row = 37; D = "'4.82"; E = "  +1.97%  " }
    @{ Row = 38; D = "'0.0296"; E = "  -3.03%  " }
    @{ Row = 39; D = "'13.18"; E = "  -6.80%  " }
    @{ Row = 40; D = $null; E = "  -5.49%  " }
    @{ Row = 41; D = "'5.86"; E = "  -0.81%  " }
    @{ Row = 42; D = "'114.23"; E = "  +17.53%  " }
    @{ Row = 43; D = "'0.202"; E = "  -5.91%  " }
    @{ Row = 44; D = "'59.93"; E = "  -2.46%  " }
    @{ Row = 45; D = "'8.80"; E = "  -3.91%  " }
    @{ Row = 46; D = $null; E = "  -2.70%  " }
    @{ Row = 47; D = $null; E = "  -0.50%  " }
    @{ Row = 48; D = $null; E = "  -4.26%  " }
    @{ Row = 49; D = "'1.16"; E = "  -1.42%  " }
    @{ Row = 50; D = "'4.24"; E = "  -13.04%  " }
    @{ Row = 51; D = "'2.25"; E = "  -0.90%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
}
